$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (ANr, KNr, Ausprägung De, Ausprägung En) to append
# starting at row 488, matching the style of existing data rows.
$rows = @(
    @("A_TYPEAREA_TRANSPORT", "K_TYPEAREA", "Verkehrsfläche", "Verkehrsfläche"),
    @("A_TYPEAREA_SETTLEMEN", "K_TYPEAREA", "Siedlungsfläche (bebaut)", "Settlement area (built on)"),
    @("A_TYPEAREA_RECREATIO", "K_TYPEAREA", "Erholungsfläche, Friedhof", "Recreation area, cemetery"),
    @("A_SERIES_FINENTRANS", "K_SERIES", "Endenergieverbrauch durch Güterbeförderung", "Final energy consumption for the transport of goods"),
    @("A_SERIES_GOODSTRANS", "K_SERIES", "Güterbeförderungsleistung", "Goods transport performance"),
    @("A_SERIES_ENCONSPTK", "K_SERIES", "Energieverbrauch je Tonnenkilometer", "Energy consumption per tonne-kilometre"),
    @("A_SERIES_PASSTRANS", "K_SERIES", "Personenbeförderungsleistung", "Passenger transport performance"),
    @("A_SERIES_ENPERPK", "K_SERIES", "Energieverbrauch je Personenkilometer", "Energy consumption per passenger-kilometre"),
    @("A_SERIES_ENINPASSTRA", "K_SERIES", "Endenergieverbrauch im Personenverkehr", "Energy consumption in passenger transport"),
    @("A_SERIES_4YAVERAGE", "K_SERIES", "Gleitender Vierjahresdurchschnitt", "Moving four-year average"),
    @("A_SERIES_YEARINQUEST", "K_SERIES", "Bezogen auf die Mittel- und Oberzentren des jeweiligen Jahres", "For the medium-sized and major cities of the year in question"),
    @("A_SERIES_YEAR2012", "K_SERIES", "Bezogen auf die Mittel- und Oberzentren des Jahres 2012", "For the medium-sized and major cities of 2012"),
    @("A_SERIES_TOTALOBJ", "K_SERIES", "Objekte insgesamt", "All objects"),
    @("A_SERIES_DIGITOBJ", "K_SERIES", "Objekte mit Digitalisat", "Objects with digitised media")
)

$startRow = 488
$srcRange = $ws.Range("A487:D487")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $destRow = $ws.Range("A" + $r + ":D" + $r)
    $srcRange.Copy($destRow)
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
    $ws.Cells.Item($r, 4).Value = $rows[$i][3]
}

# Widen column C to fit the newly added (longer) German descriptions.
$ws.Columns.Item(3).ColumnWidth = 71
